$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case data added to row 2 (QA-3 test cases)
$ws.Range("C2").Value = "RbcLd_2704732"
$ws.Range("D2").Value = "USDNexnDLd"
$ws.Range("I2").Value = "V_RYORJHML"
$ws.Range("K2").Value = "jDIubtQImU"
$ws.Range("M2").Value = "jDIubtQImU"
$ws.Range("N2").Value = "UCN 10417"
